$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in values for R2, S2, T2, U2 (T2 keeps its existing style, R2/S2/U2 have no explicit style)
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 5

# Row 3: R3, S3, T3 already styled; just set their values. U3 is a brand-new plain cell.
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = 5

# Row 7: R7, S7, T7 get values; U7 is new.
$ws.Range("R7").Value = 5
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 5
$ws.Range("U7").Value = 5

# Row 8: only U8 is new.
$ws.Range("U8").Value = 5

# Row 10: only U10 is new.
$ws.Range("U10").Value = 5

# Row 17: T17 gets a value; U17 is new.
$ws.Range("T17").Value = 5
$ws.Range("U17").Value = 5

# Row 19: U19 is new and carries the same direct formatting (style index 12) as V4,
# so clone the format from V4 before setting the value.
$ws.Range("V4").Copy()
$ws.Range("U19").PasteSpecial(-4122)
$ws.Range("U19").Value = 5

# Row 20: R20, S20, T20 get values; U20 is new.
$ws.Range("R20").Value = 5
$ws.Range("S20").Value = 5
$ws.Range("T20").Value = 5
$ws.Range("U20").Value = 5

# Row 22: R22, S22, T22 get values; U22 is new.
$ws.Range("R22").Value = 5
$ws.Range("S22").Value = 5
$ws.Range("T22").Value = 5
$ws.Range("U22").Value = 5

# Update the selection / active cell shown when the sheet is reopened.
$ws.Range("T3").Select()
